# Update rows 2-16 of Sheet1 (S100a9 -> Tlr4 natmi edge table) with recomputed
# ligand/receptor specificity values, and add a new "ECs" sending-cluster block
# (rows 12-16), per Dr Hou's revised NATMI analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A:T for each row -> Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, then the 16 numeric NATMI metric columns (E:T).
$rows = @{
    2 = @("M1", "S100a9", "Tlr4", "ECs", 2, 0.6666666666666666, 0.1102606666666667, 0.330782, 0.001202219334792345, 0.001202219334792345, 3, 1, 7.908652666666666, 23.725958, 0.1059802478036631, 0.1086148104297394, 0.8720133154617777, 7.848119839156, 0.0001274115030156478, 0.000130578825143438)
    3 = @("M1", "S100a9", "Tlr4", "FAPs", 2, 0.6666666666666666, 0.1102606666666667, 0.330782, 0.001202219334792345, 0.001202219334792345, 3, 1, 10.76975733333333, 32.309272, 0.1443206066922968, 0.1479082721718924, 1.187480623411556, 10.687325610704, 0.0001735050237744408, 0.0001778181845807776)
    4 = @("M1", "S100a9", "Tlr4", "M1", 2, 0.6666666666666666, 0.1102606666666667, 0.330782, 0.001202219334792345, 0.001202219334792345, 3, 1, 19.37830233333333, 58.134907, 0.2596797924831069, 0.2661351715768666, 2.136664534141556, 19.229980807274, 0.000312192067378055, 0.0003199528489379872)
    5 = @("M1", "S100a9", "Tlr4", "M2", 2, 0.6666666666666666, 0.1102606666666667, 0.330782, 0.001202219334792345, 0.001202219334792345, 3, 1, 31.136893, 93.41067899999999, 0.4172513037378061, 0.4276237525034073, 3.433174580108667, 30.898571220978, 0.0005016275848209041, 0.0005140975432760529)
    6 = @("M1", "S100a9", "Tlr4", "sCs", 2, 0.6666666666666666, 0.1102606666666667, 0.330782, 0.001202219334792345, 0.001202219334792345, 2, 1, 5.430231, 10.860462, 0.07276804928312694, 0.0497179933180944, 0.598740890214, 3.592445341284, 0.000087483155803297471455465262, 0.000059771932854089717952576766)
    7 = @("M2", "S100a9", "Tlr4", "ECs", 3, 1, 6.461386, 19.384158, 0.07045126257253938, 0.07045126257253938, 3, 1, 7.908652666666666, 23.725958, 0.1059802478036631, 0.1086148104297394, 51.10085761926266, 459.907718573364, 0.007466442265518659, 0.007652050528852158)
    8 = @("M2", "S100a9", "Tlr4", "FAPs", 3, 1, 6.461386, 19.384158, 0.07045126257253938, 0.07045126257253938, 3, 1, 10.76975733333333, 32.309272, 0.1443206066922968, 0.1479082721718924, 69.58755925699732, 626.2880333129759, 0.01016756895670719, 0.01042032451943261)
    9 = @("M2", "S100a9", "Tlr4", "M1", 3, 1, 6.461386, 19.384158, 0.07045126257253938, 0.07045126257253938, 3, 1, 19.37830233333333, 58.134907, 0.2596797924831069, 0.2661351715768666, 125.2106914003673, 1126.896222603306, 0.0182947692450099, 0.01874955885254964)
    10 = @("M2", "S100a9", "Tlr4", "M2", 3, 1, 6.461386, 19.384158, 0.07045126257253938, 0.07045126257253938, 3, 1, 31.136893, 93.41067899999999, 0.4172513037378061, 0.4276237525034073, 201.187484513698, 1810.687360623282, 0.02939588115836655, 0.03012663326987214)
    11 = @("M2", "S100a9", "Tlr4", "sCs", 3, 1, 6.461386, 19.384158, 0.07045126257253938, 0.07045126257253938, 2, 1, 5.430231, 10.860462, 0.07276804928312694, 0.0497179933180944, 35.086818560166, 210.520911360996, 0.005126600946937062, 0.003502695401832827)
    12 = @("ECs", "S100a9", "Tlr4", "ECs", 3, 1, 85.14262166666667, 255.427865, 0.9283465180926683, 0.9283465180926682, 3, 1, 7.908652666666666, 23.725958, 0.1059802478036631, 0.1086148104297394, 673.3634218910744, 6060.27079701967, 0.0983863940351288, 0.1008321810757438)
    13 = @("ECs", "S100a9", "Tlr4", "FAPs", 3, 1, 85.14262166666667, 255.427865, 0.9283465180926683, 0.9283465180926682, 3, 1, 10.76975733333333, 32.309272, 0.1443206066922968, 0.1479082721718924, 916.9653740738089, 8252.688366664281, 0.1339795327118152, 0.137310129467879)
    14 = @("ECs", "S100a9", "Tlr4", "M1", 3, 1, 85.14262166666667, 255.427865, 0.9283465180926683, 0.9283465180926682, 3, 1, 19.37830233333333, 58.134907, 0.2596797924831069, 0.2661351715768666, 1649.919464109284, 14849.27517698355, 0.2410728311707189, 0.2470656598753789)
    15 = @("ECs", "S100a9", "Tlr4", "M2", 3, 1, 85.14262166666667, 255.427865, 0.9283465180926683, 0.9283465180926682, 3, 1, 31.136893, 93.41067899999999, 0.4172513037378061, 0.4276237525034073, 2651.076700574481, 23859.69030517033, 0.3873537949946186, 0.3969830216902591)
    16 = @("ECs", "S100a9", "Tlr4", "sCs", 3, 1, 85.14262166666667, 255.427865, 0.9283465180926683, 0.9283465180926682, 2, 1, 5.430231, 10.860462, 0.07276804928312694, 0.0497179933180944, 462.344103595605, 2774.06462157363, 0.06755396518038659, 0.04615552598340748)
}

foreach ($r in $rows.Keys) {
    $rowData = $rows[$r]
    for ($i = 0; $i -lt $rowData.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $rowData[$i]
    }
}
